# Apply the "RotJ - clean up folder, replace my wip with Sonikkustars 4 frame
# faster version" edit to the FrameCounts sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab.
$ws.Name = "FrameCounts"

# Section header (row 5): "Circuits" -> "Level 1"
$ws.Range("A5").Value = "Level 1"

# Header row: columns B/C get new labels ("Mine" / "Andymac"); A1/D1 unchanged.
$ws.Range("B1").Value = "Mine"
$ws.Range("C1").Value = "Andymac"

# Data rows 6-12: new labels + new frame numbers.
$ws.Range("A6").Value = "Batman appears"
$ws.Range("B6").Value = 450
$ws.Range("C6").Value = 450

$ws.Range("A7").Value = "X = 210"
$ws.Range("B7").Value = 570
$ws.Range("C7").Value = 622

$ws.Range("A8").Value = "X = 579"
$ws.Range("B8").Value = 690
$ws.Range("C8").Value = 741

$ws.Range("A9").Value = "X = 901"
$ws.Range("B9").Value = 844
$ws.Range("C9").Value = 900

$ws.Range("A10").Value = "X = 1129"
$ws.Range("B10").Value = 988
$ws.Range("C10").Value = 1040

$ws.Range("A11").Value = "X = 1252"
$ws.Range("B11").Value = 1281
$ws.Range("C11").Value = 1330

$ws.Range("A12").Value = "Screen 2"
$ws.Range("B12").Value = 1666
$ws.Range("C12").Value = 1677

# Rows 13-19 no longer carry any data (A/B/C cleared, D keeps its formula and
# naturally recomputes to 0 since B is blank).
$ws.Range("A13:C19").Clear()

# Column C needs a slightly wider custom width now that it holds longer values.
$ws.Columns("C").ColumnWidth = 9.5

# Restore the view: scroll back up (frozen pane resets topLeftCell to A2) and
# move the active selection to B13.
[void]$ws.Range("B13").Select()
